# DPLKKPS003-001: add "Verifikasi Register Deposit (Bulk)" (row 3) and
# "Verifikasi Deposit (Bulk)" (row 4) scripted scenarios, plus two new
# STATUS_VERIFIKASI / KETERANGAN_VERIFIKASI columns (U/V).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 2: RUN marker moves down to row 4, so A2 is cleared out.
# ---------------------------------------------------------------------
$ws.Range("A2").Clear()

# ---------------------------------------------------------------------
# 2. New header cells for the verification columns.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("U1").PasteSpecial(-4122)
$ws.Range("V1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Row heights for the two new data rows.
# ---------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 105
$ws.Rows.Item(4).RowHeight = 105

# ---------------------------------------------------------------------
# 4. Formats for the newly-populated cells in row 3 (copy style from the
#    matching column in row 2, which already carries the right xf).
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("M3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5. Formats for the newly-populated cells in row 4.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("L4").PasteSpecial(-4122)

$ws.Range("N2").Copy()
$ws.Range("N4").PasteSpecial(-4122)

$ws.Range("O2").Copy()
$ws.Range("O4").PasteSpecial(-4122)
$ws.Range("Q4").PasteSpecial(-4122)

$ws.Range("P2").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("V4").PasteSpecial(-4122)

$ws.Range("R2").Copy()
$ws.Range("R4").PasteSpecial(-4122)

# Row 4 no longer carries the STATUS_REGISTER / KETERANGAN_REGISTER cells.
$ws.Range("S4").Clear()
$ws.Range("T4").Clear()

# ---------------------------------------------------------------------
# 6. Cell values. Written in the same order the source workbook used so
#    the shared-string table grows the same way.
# ---------------------------------------------------------------------
$ws.Range("E3").Value = "Verifikasi Register Deposit (Bulk)"
$ws.Range("I4").Value = "Penyelia Settlement"
$ws.Range("K4").Value = "Proses"
$ws.Range("E4").Value = "Verifikasi Deposit (Bulk)"
$ws.Range("S3").Value = "1 : Lanjutkan ke Verifikasi"
$ws.Range("N3").Value = "000007947"
$ws.Range("U4").Value = "1 : Setuju"
$ws.Range("V4").Value = "Setuju Verifikasi"
$ws.Range("F4").Value = "Username : 31224;`nPassword : bni1234; `nKode Perusahaan : 000007947;`nStatus Register : 1 : Setuju;`nKeterangan Upload : Setuju Verifikasi`n"
$ws.Range("U1").Value = "STATUS_VERIFIKASI"
$ws.Range("V1").Value = "KETERANGAN_VERIFIKASI"
$ws.Range("T3").Value = "KEP.TRX.086/24"
$ws.Range("F3").Value = "Username : 33028;`nPassword : bni1234; `nKode Perusahaan : 000007947;`nStatus Register : 1 : Lanjutkan ke Verifikasi;`nKeterangan Register : KEP.TRX.086/24"

# Remaining row 3 values (re-use existing shared strings).
$ws.Range("B3").Value = "DPLKKPS003-001"
$ws.Range("C3").Value = "Normal - Kepesertaan - Transaksi"
$ws.Range("D3").Value = "Kepesertaan - Transaksi"
$ws.Range("G3").Value = 33028
$ws.Range("H3").Value = "bni1234"
$ws.Range("I3").Value = "Asistent Settlement"
$ws.Range("J3").Value = "Kepesertaan"
$ws.Range("K3").Value = "Transaksi"
$ws.Range("L3").Value = "Deposit"
$ws.Range("M3").Value = "Register Deposit (Bulk)"

# Remaining row 4 values.
$ws.Range("A4").Value = "RUN"
$ws.Range("B4").Value = "DPLKKPS003-001"
$ws.Range("C4").Value = "Normal - Kepesertaan - Transaksi"
$ws.Range("D4").Value = "Kepesertaan - Transaksi"
$ws.Range("G4").Value = 30711
$ws.Range("H4").Value = "bni1234"
$ws.Range("J4").Value = "Kepesertaan"
$ws.Range("L4").Value = "Deposit"
$ws.Range("M4").Value = "Verifikasi Deposit (Bulk)"
$ws.Range("N4").Value = "000007947"

# ---------------------------------------------------------------------
# 7. Column widths for the new U/V columns (best-fit on the new content).
# ---------------------------------------------------------------------
$ws.Columns.Item(21).ColumnWidth = 17.59
$ws.Columns.Item(22).ColumnWidth = 23.17

# ---------------------------------------------------------------------
# 8. Restore view/selection to what the scripted session left behind.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("H4").Select()
